$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the header row (row 1) to insert new "... sentiment" columns
# immediately after each of the existing aspect columns (plot/story, acting,
# direction, cinemetography, editing, music, character/world development).
$headers = @{
    "A1" = "sentence"
    "B1" = "part of article"
    "C1" = "sentiment (pos = 1, neg = -1, neutral = 0)"
    "D1" = "plot/story"
    "E1" = "plot/story sentiment"
    "F1" = "acting"
    "G1" = "acting sentiment"
    "H1" = "direction"
    "I1" = "direction sentiment"
    "J1" = "cinemetography"
    "K1" = "cinemetography sentiment"
    "L1" = "editing"
    "M1" = "editing sentiment"
    "N1" = "music"
    "O1" = "music sentiment"
    "P1" = "character/world development"
    "Q1" = "character/world development sentiment"
    "R1" = "#positve sentiment words "
    "S1" = "#negative sentiment words "
    "T1" = "confusing sentiment words list "
    "U1" = "name of the movie "
    "V1" = "method used"
    "W1" = "Website"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Update the view: select H3 instead of H18 (and scroll the window so
# the selection is visible, matching the new topLeftCell="H1").
$excel.Goto($ws.Range("H3"), $true)
